$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.531.32"
$ws.Range("E2").Value = "  +3.04%  "
$ws.Range("D3").Value = "3.068.05"
$ws.Range("E3").Value = "  +3.26%  "
$ws.Range("E4").Value = "  +0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "518.23"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +3.28%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "141.14"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +3.21%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +1.92%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "7.24"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("E10").Value = "  +1.52%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.375"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +3.54%  "
$ws.Range("D12").Value = "3.598.21"
$ws.Range("E12").Value = "  +3.34%  "
$ws.Range("E13").Value = "  +3.27%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "25.56"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.19%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.0000162"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").Value = "57.614.46"
$ws.Range("E16").Value = "  +3.09%  "
$ws.Range("D17").Value = "3.067.20"
$ws.Range("E17").Value = "  +3.37%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "6.10"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +2.12%  "
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("E20").Value = "  +1.89%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "331.65"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("E23").Value = "  +1.59%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "65.70"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.19%  "
$ws.Range("E25").Value = "  +4.84%  "
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").Value = "  +1.47%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "6.37"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("E29").Value = "  +3.93%  "
$ws.Range("E30").Value = "  +3.25%  "
$ws.Range("E31").Value = "  +4.32%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "20.73"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +3.05%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "154.74"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.89%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.49"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.49%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "27.05"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +6.43%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "5.91"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +4.30%  "
$ws.Range("E37").Value = "  +2.90%  "
$ws.Range("E38").Value = "  +2.76%  "
$ws.Range("D39").Value = "3.109.39"
$ws.Range("E39").Value = "  +3.43%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.92"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +4.44%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "36.76"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.10%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.04%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.654"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("D44").Value = "2.262.32"
$ws.Range("E44").Value = "  +5.13%  "
$ws.Range("E45").Value = "  +10.78%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "20.89"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +7.41%  "
$ws.Range("E47").Value = "  +1.56%  "
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("E49").Value = "  +1.10%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "262.69"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +15.83%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.706"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +5.88%  "
